# "ajout de la selection de stat perso"
# Adds 4 new rows (hero stats) to the "donnees" sheet and switches the
# active sheet / selection to that sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("donnees")

# --- new hero-stat rows -----------------------------------------------
$ws.Range("A4").Value = "atk hero"
$ws.Range("B4").Value = 10000

$ws.Range("A5").Value = "def hero"
$ws.Range("B5").Value = 0

$ws.Range("A6").Value = "pv hero"
$ws.Range("B6").Value = 1001

$ws.Range("A7").Value = "% augmentation stats par level"

# B7 must be stored as the literal text "10.0" (not the number 10).
# Writing the string straight into .Value auto-coerces a numeric-looking
# string into a number, and flipping NumberFormat to force text leaves a
# permanent (unused) style behind. Instead, build the text via a formula
# in a scratch cell, then copy/paste-special just the resulting value —
# that yields a plain shared-string cell with no style side effects.
$scratch = $ws.Range("D1")
$scratch.Formula = '="10.0"'
$scratch.Copy()
$ws.Range("B7").PasteSpecial(-4163) # xlPasteValues
$scratch.ClearContents()

# --- selection / active sheet ------------------------------------------
# Make "donnees" the active sheet (moves tabSelected from "pokemon" to
# "donnees" and updates workbook.xml's activeTab automatically), and
# select the freshly added block.
$ws.Activate()
$ws.Range("A4:B7").Select()
